$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GKW"
$ws.Range("B2").Value = "June 27, 2023"
$ws.Range("E2").Value = 61753

$ws.Range("A3").Value = "MAF"
$ws.Range("B3").Value = "June 27, 2023"
$ws.Range("E3").Value = 4429

$ws.Range("A4").Value = "SAEM"
$ws.Range("B4").Value = "June 27, 2023"

$ws.Range("A5").Value = "CA"
$ws.Range("B5").Value = "June 27, 2023"
$ws.Range("E5").Value = 399132
$ws.Range("F5").Value = 115

$ws.Range("A6").Value = "MV"
$ws.Range("B6").Value = "June 27, 2023"
$ws.Range("E6").Value = 4754690

$ws.Range("A7").Value = "MC"
$ws.Range("B7").Value = "June 27, 2023"
$ws.Range("E7").Value = 43158
$ws.Range("F7").Value = 92

$ws.Range("A8").Value = "MD"
$ws.Range("B8").Value = "June 27, 2023"
$ws.Range("E8").Value = 322274

$ws.Range("A9").Value = "TM"
$ws.Range("B9").Value = "June 27, 2023"
$ws.Range("E9").Value = 2412814

$ws.Range("A10").Value = "XM"
$ws.Range("B10").Value = "June 27, 2023"
$ws.Range("E10").Value = 42795
$ws.Range("F10").Value = 95

$ws.Range("A11").Value = "SNLS"
$ws.Range("B11").Value = "June 27, 2023"
$ws.Range("E11").Value = 24969
$ws.Range("F11").Value = 62

$ws.Range("A12").Value = "MSA"
$ws.Range("B12").Value = "June 27, 2023"
$ws.Range("E12").Value = 2010530
$ws.Range("F12").Value = 13

$ws.Range("A13").Value = "MSP"
$ws.Range("B13").Value = "June 27, 2023"
$ws.Range("E13").Value = 1047308

$ws.Range("A14").Value = "VMF"
$ws.Range("B14").Value = "June 27, 2023"
$ws.Range("E14").Value = 258086
$ws.Range("F14").Value = 110

$ws.Range("A15").Value = "MCL"
$ws.Range("B15").Value = "June 27, 2023"
$ws.Range("E15").Value = 116054
$ws.Range("F15").Value = 193

$ws.Range("A16").Value = "VM"
$ws.Range("B16").Value = "June 27, 2023"
$ws.Range("E16").Value = 77100
$ws.Range("F16").Value = 157

$ws.Range("A17").Value = "SU"
$ws.Range("B17").Value = "June 27, 2023"

$ws.Range("A18").Value = "STVA"
$ws.Range("B18").Value = "June 27, 2023"
$ws.Range("E18").Value = 411772
$ws.Range("F18").Value = 50

$ws.Range("A19").Value = "ATT"
$ws.Range("B19").Value = "June 27, 2023"
$ws.Range("E19").Value = 3991306
$ws.Range("F19").Value = 15

$ws.Range("A20").Value = "SSLG"
$ws.Range("B20").Value = "June 27, 2023"
$ws.Range("E20").Value = 4150

$ws.Range("A21").Value = "MDCM"
$ws.Range("B21").Value = "June 27, 2023"
$ws.Range("E21").Value = 6605

$ws.Range("A22").Value = "MM"
$ws.Range("B22").Value = "June 27, 2023"
$ws.Range("E22").Value = 3129
$ws.Range("F22").Value = 61

$ws.Range("A23").Value = "MFR"
$ws.Range("B23").Value = "June 27, 2023"
$ws.Range("E23").Value = 45222

$ws.Range("A24").Value = "XF"
$ws.Range("B24").Value = "June 27, 2023"
$ws.Range("E24").Value = 877238
$ws.Range("F24").Value = 7

$ws.Range("A25").Value = "GFBR"
$ws.Range("B25").Value = "June 27, 2023"
$ws.Range("E25").Value = 295

$ws.Range("A26").Value = "MVIA"
$ws.Range("B26").Value = "June 27, 2023"
$ws.Range("E26").Value = 1622

$ws.Range("A27").Value = "ARM"
$ws.Range("B27").Value = "June 27, 2023"

$ws.Range("A28").Value = "ASTRCN"
$ws.Range("B28").Value = "June 27, 2023"

$ws.Range("A29").Value = "HUGH"
$ws.Range("B29").Value = "June 27, 2023"

$ws.Range("A30").Value = "HTMYA"
$ws.Range("B30").Value = "June 27, 2023"

$ws.Range("A31").Value = "MIDCO"
$ws.Range("B31").Value = "June 27, 2023"

$ws.Range("A32").Value = "OPTS"
$ws.Range("B32").Value = "June 27, 2023"

$ws.Range("A33").Value = "USCELL"
$ws.Range("B33").Value = "June 27, 2023"
$ws.Range("E33").Value = 28374

$ws.Range("A34").Value = "SEC"
$ws.Range("B34").Value = "June 27, 2023"

$ws.Range("A35").Value = "OPTTV"
$ws.Range("B35").Value = "June 27, 2023"
$ws.Range("E35").Value = 13443

$ws.Range("A36").Value = "BRE"
$ws.Range("B36").Value = "June 27, 2023"

$ws.Range("A37").Value = "BLUER"
$ws.Range("B37").Value = "June 27, 2023"
$ws.Range("E37").Value = 3300

$ws.Range("A38").Value = "BUCK"
$ws.Range("B38").Value = "June 27, 2023"
